$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'30.497.91"
$ws.Range("E2").Value = '  +0.36%  '

$ws.Range("D3").Value = "'2.105.25"
$ws.Range("E3").Value = '  +4.52%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").Value = "'329.98"
$ws.Range("E5").Value = '  +1.44%  '

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = '  +0.16%  '

$ws.Range("D7").Value = "'0.5264"
$ws.Range("E7").Value = '  +2.68%  '

$ws.Range("D8").Value = "'0.4396"
$ws.Range("E8").Value = '  +3.04%  '

$ws.Range("D9").Value = "'0.08872"
$ws.Range("E9").Value = '  +1.52%  '

$ws.Range("D10").Value = "'47.59"
$ws.Range("E10").Value = '  +9.70%  '

$ws.Range("D11").Value = "'1.165"
$ws.Range("E11").Value = '  +2.53%  '

$ws.Range("D12").Value = "'24.64"
$ws.Range("E12").Value = '  -0.14%  '

$ws.Range("D13").Value = "'2.114.63"
$ws.Range("E13").Value = '  +5.04%  '

$ws.Range("E14").Value = '  +2.10%  '

$ws.Range("E15").Value = '  +4.11%  '

$ws.Range("D16").Value = "'96.45"
$ws.Range("E16").Value = '  +2.41%  '

$ws.Range("E17").Value = '  +0.32%  '

$ws.Range("D19").Value = "'0.06645"
$ws.Range("E19").Value = '  +1.86%  '

$ws.Range("E20").Value = '  +0.75%  '

$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = '  +0.16%  '

$ws.Range("E22").Value = '  +1.41%  '

$ws.Range("D23").Value = "'30.552.72"
$ws.Range("E23").Value = '  +0.35%  '

$ws.Range("E24").Value = '  +3.99%  '

$ws.Range("E25").Value = '  +3.96%  '

$ws.Range("D26").Value = "'2.357.42"
$ws.Range("E26").Value = '  +4.79%  '

$ws.Range("D27").Value = "'22.44"
$ws.Range("E27").Value = '  -0.10%  '

$ws.Range("D28").Value = "'2.602"
$ws.Range("E28").Value = '  +6.68%  '

$ws.Range("D29").Value = "'161.70"
$ws.Range("E29").Value = '  -0.54%  '

$ws.Range("D30").Value = "'132.75"
$ws.Range("E30").Value = '  +1.24%  '

$ws.Range("D31").Value = "'1.211"
$ws.Range("E31").Value = '  +5.46%  '

$ws.Range("E32").Value = '  +1.97%  '

$ws.Range("D33").Value = "'1.681"
$ws.Range("E33").Value = '  +22.33%  '

$ws.Range("D34").Value = "'6.227"
$ws.Range("E34").Value = '  +2.06%  '

$ws.Range("D35").Value = "'3.934"
$ws.Range("E35").Value = '  +2.72%  '

$ws.Range("D36").Value = "'10.15"
$ws.Range("E36").Value = '  +11.12%  '

$ws.Range("D37").Value = "'0.02584"
$ws.Range("E37").Value = '  +2.08%  '

$ws.Range("D38").Value = "'5.502"
$ws.Range("E38").Value = '  +0.73%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = "'0.06707"
$ws.Range("E39").Value = '  +0.65%  '

$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = "'12.75"
$ws.Range("E40").Value = '  +2.83%  '

$ws.Range("D41").Value = "'0.2284"
$ws.Range("E41").Value = '  +3.40%  '

$ws.Range("D42").Value = "'0.6861"
$ws.Range("E42").Value = '  +2.85%  '

$ws.Range("D43").Value = "'1.271"
$ws.Range("E43").Value = '  +2.63%  '

$ws.Range("E44").Value = '  +0.12%  '

$ws.Range("D45").Value = "'14.08"
$ws.Range("E45").Value = '  +3.15%  '

$ws.Range("D46").Value = "'0.6394"
$ws.Range("E46").Value = '  +3.60%  '

$ws.Range("D47").Value = "'2.216"
$ws.Range("E47").Value = '  +0.78%  '

$ws.Range("D48").Value = "'3.631"
$ws.Range("E48").Value = '  -1.11%  '

$ws.Range("E49").Value = '  -0.72%  '

$ws.Range("D50").Value = "'1.211"
$ws.Range("E50").Value = '  +9.50%  '

$ws.Range("D51").Value = "'82.31"
$ws.Range("E51").Value = '  +1.71%  '
